# Fruta / hortaliza, semanal
# Inserts a new weekly price record at row 188 (pushing the existing
# rows 188-213 down to 189-214) for "Agrícola del Norte S.A. de Arica - Plátano".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before the current row 188; Excel shifts rows
# 188-213 down to 189-214 automatically (formats copied from the row above).
$ws.Rows(188).Insert()

# Populate the new row 188 with the new weekly observation.
$ws.Cells.Item(188, 1).Value  = 1
$ws.Cells.Item(188, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(188, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(188, 4).Value  = 44644
$ws.Cells.Item(188, 5).Value  = 15
$ws.Cells.Item(188, 6).Value  = "Fruta"
$ws.Cells.Item(188, 7).Value  = 100108
$ws.Cells.Item(188, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(188, 9).Value  = 100108006
$ws.Cells.Item(188, 10).Value = "Plátano"
$ws.Cells.Item(188, 11).Value = "Sin especificar"
$ws.Cells.Item(188, 12).Value = "Pintón"
$ws.Cells.Item(188, 13).Value = 120
$ws.Cells.Item(188, 14).Value = 19000
$ws.Cells.Item(188, 15).Value = 20000
$ws.Cells.Item(188, 16).Value = 19500
$ws.Cells.Item(188, 17).Value = "`$/caja 20 kilos"
$ws.Cells.Item(188, 18).Value = "Ecuador"
$ws.Cells.Item(188, 19).Value = 975
$ws.Cells.Item(188, 20).Value = 20
